$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 148
$ws.Range("I2").Value = 390
$ws.Range("J2").Value = 1668
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 459
$ws.Range("M2").Value = 21
$ws.Range("N2").Value = 301
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 200
$ws.Range("T2").Value = 277
$ws.Range("U2").Value = 23
$ws.Range("V2").Value = 2530
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2473
$ws.Range("Y2").Value = 0
$ws.Range("AA2").Value = 19
